# Weekly update: shift the "Fecha" (and, where applicable, "Origen") values
# of the Perejil / Vega Monumental Concepción data block (rows 132-173) up
# by one reporting period, refresh the newest period's prices (rows 132-133)
# with the latest figures, and append the previous oldest period's record
# (rows 174-175) so the rolling weekly window stays populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New values for the most recent period (rows 132-133) ---
$ws.Range("D132").Value = 44841
$ws.Range("K132").Value = 700
$ws.Range("L132").Value = 800
$ws.Range("M132").Value = 750
$ws.Range("P132").Value = 750

$ws.Range("D133").Value = 44841
$ws.Range("K133").Value = 600
$ws.Range("L133").Value = 600
$ws.Range("M133").Value = 600
$ws.Range("P133").Value = 600

# --- Dates (and Volumen J where it moved with the date) shift up by one
#     reporting period for rows 134-173 ---
$ws.Range("D134").Value = 44435
$ws.Range("J134").Value = 200
$ws.Range("D135").Value = 44435
$ws.Range("J135").Value = 100

$ws.Range("D136").Value = 44442
$ws.Range("J136").Value = 300
$ws.Range("D137").Value = 44442
$ws.Range("J137").Value = 150

$ws.Range("D138").Value = 44336
$ws.Range("D139").Value = 44336

$ws.Range("D140").Value = 44252
$ws.Range("D141").Value = 44252

$ws.Range("D142").Value = 44694
$ws.Range("D143").Value = 44694

$ws.Range("D144").Value = 44405
$ws.Range("D145").Value = 44405

$ws.Range("D146").Value = 44679
$ws.Range("D147").Value = 44679

$ws.Range("D148").Value = 44231
$ws.Range("D149").Value = 44231

$ws.Range("D150").Value = 44334
$ws.Range("D151").Value = 44334

$ws.Range("D152").Value = 44194
$ws.Range("D153").Value = 44194

$ws.Range("D154").Value = 44330
$ws.Range("D155").Value = 44330

$ws.Range("D156").Value = 44274
$ws.Range("D157").Value = 44274

$ws.Range("D158").Value = 44391
$ws.Range("D159").Value = 44391

$ws.Range("D160").Value = 44433
$ws.Range("D161").Value = 44433

$ws.Range("D162").Value = 44203
$ws.Range("D163").Value = 44203

$ws.Range("D164").Value = 44355
$ws.Range("O164").Value = "Región de Ñuble"
$ws.Range("D165").Value = 44355
$ws.Range("O165").Value = "Región de Ñuble"

$ws.Range("D166").Value = 44565
$ws.Range("O166").Value = "Región Metropolitana"
$ws.Range("D167").Value = 44565
$ws.Range("O167").Value = "Región Metropolitana"

$ws.Range("D168").Value = 44187
$ws.Range("D169").Value = 44187

$ws.Range("D170").Value = 44553
$ws.Range("D171").Value = 44553

$ws.Range("D172").Value = 44292
$ws.Range("D173").Value = 44292

# --- Append the record that fell off the top of the window (previously
#     rows 132-133) as new rows 174-175, dated with the next period ---
$ws.Range("A174").Value = 11
$ws.Range("B174").Value = "Vega Monumental Concepción"
$ws.Range("C174").Value = "Bíobío"
$ws.Range("D174").Value = 44453
$ws.Range("E174").Value = 8
$ws.Range("F174").Value = 100112044
$ws.Range("G174").Value = "Perejil"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 200
$ws.Range("K174").Value = 600
$ws.Range("L174").Value = 700
$ws.Range("M174").Value = 650
$ws.Range("N174").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O174").Value = "Región de Ñuble"
$ws.Range("P174").Value = 650
$ws.Range("Q174").Value = 1
$ws.Range("R174").Value = "Hortaliza"

$ws.Range("A175").Value = 11
$ws.Range("B175").Value = "Vega Monumental Concepción"
$ws.Range("C175").Value = "Bíobío"
$ws.Range("D175").Value = 44453
$ws.Range("E175").Value = 8
$ws.Range("F175").Value = 100112044
$ws.Range("G175").Value = "Perejil"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Segunda"
$ws.Range("J175").Value = 100
$ws.Range("K175").Value = 500
$ws.Range("L175").Value = 500
$ws.Range("M175").Value = 500
$ws.Range("N175").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O175").Value = "Región de Ñuble"
$ws.Range("P175").Value = 500
$ws.Range("Q175").Value = 1
$ws.Range("R175").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D174").NumberFormat = $ws.Range("D173").NumberFormat
$ws.Range("D175").NumberFormat = $ws.Range("D173").NumberFormat
